# Add capital costs for petrol and LPG buses (PR #174).
# Replaces the old trailing blank/placeholder row (previously row 61, which
# only carried a leftover cell style) with two new data rows describing
# "Assumed same as Diesel ICE" bus capital-cost entries for Petrol ICE and
# LPG technologies, each costed at 400000 with a 0.86 servicing/tyres ratio
# (mirroring the existing Bus/Diesel ICE row directly above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 61 held a single formatted-but-empty cell (G61). Clear its
# formatting before repurposing the row for real data so no stray style
# lingers on the (now populated) G61 cell.
$ws.Range("G61").ClearFormats()

# --- Row 61: Bus / Petrol ICE ---------------------------------------------
$ws.Range("B61").Value = "Assumed same as Diesel ICE"
$ws.Range("C61").Value = 2025
$ws.Range("D61").Value = "Bus"
$ws.Range("E61").Value = "Bus"
$ws.Range("F61").Value = "Petrol ICE"
$ws.Range("G61").Value = 400000
# K column uses the same "0.00" number style as the other literal (non
# formula) ratio cells further up the table (e.g. K56) - copy that format
# across before writing the value.
$ws.Range("K56").Copy()
$ws.Range("K61").PasteSpecial(-4122)
$ws.Range("K61").Value = 0.86

# --- Row 62: Bus / LPG -----------------------------------------------------
$ws.Range("B62").Value = "Assumed same as Diesel ICE"
$ws.Range("C62").Value = 2025
$ws.Range("D62").Value = "Bus"
$ws.Range("E62").Value = "Bus"
$ws.Range("F62").Value = "LPG"
$ws.Range("G62").Value = 400000
$ws.Range("K56").Copy()
$ws.Range("K62").PasteSpecial(-4122)
$ws.Range("K62").Value = 0.86

$excel.CutCopyMode = $false

# Re-apply the autofilter over the now-larger table (A1:K59 -> A1:K63) and
# keep the workbook-level hidden _FilterDatabase defined name in sync with
# it (Excel keeps these two in lock-step whenever the filter range grows).
$ws.AutoFilterMode = $false
$ws.Range("A1:K63").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "AG_costs!_FilterDatabase") {
        $n.RefersTo = "=AG_costs!`$A`$1:`$K`$63"
    }
}

# Match the author's final on-screen selection from the saved file.
$ws.Range("H66").Select()
